$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, shifting the existing data (rows 56-134)
# down to rows 57-135.
$ws.Range("A56").EntireRow.Insert()

# Populate the newly inserted row 56 with the new record.
$ws.Range("A56").Value = 11
$ws.Range("B56").Value = "Vega Monumental Concepción"
$ws.Range("C56").Value = "Bíobío"
$ws.Range("D56").Value = 44665
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100109
$ws.Range("H56").Value = "Uva"
$ws.Range("I56").Value = 100109001
$ws.Range("J56").Value = "Uva"
$ws.Range("K56").Value = "Red Globe"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 170
$ws.Range("N56").Value = 8000
$ws.Range("O56").Value = 8500
$ws.Range("P56").Value = 8265
$ws.Range("Q56").Value = "`$/bandeja 18 kilos"
$ws.Range("R56").Value = "Región de O'Higgins"
$ws.Range("S56").Value = 459
$ws.Range("T56").Value = 18
